# Fix: waiting-display / per-toko transaction file writer was leaving a
# stale transaction row (and the header columns in the wrong order) in
# dataBill.xlsx. Reset the sheet back to a clean header-only state with
# the columns the dashboard code actually expects: toko, menu, qty, harga.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stale "last transaction" row that was left behind.
$ws.Rows("2:2").Delete()

# Re-write the header row with the corrected column names/order.
$ws.Range("A1").Value = "toko"
$ws.Range("B1").Value = "menu"
$ws.Range("C1").Value = "qty"
$ws.Range("D1").Value = "harga"

# Column B ("menu") now needs room for the longer menu names, matching
# the width the fixed writer lays the sheet out with. Columns D/E keep
# their existing saved widths (harga / qty helper column) unchanged.
$ws.Columns("B").ColumnWidth = 19.45

# Selection left on the (now empty) next row, matching a row-delete.
$ws.Range("A2:XFD2").Select() | Out-Null
